# Generate Report for Handback
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# The "Status" text used across sheets changes from "Ready for handoff"
# to "Handed back: in sync with en-US"
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# zh-cn sheet: Status, Latest Handback DateTime, Error Detail
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("K2").Value = "2016-09-05 02:55:30"
$zhcn.Range("P2").Value = ""

# de-de sheet: Status, Latest Handback DateTime, Error Detail
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("K2").Value = "2016-09-05 02:55:38"
$dede.Range("P2").Value = ""

# Column width adjustments
$overview.Range("E:E").ColumnWidth = 29.9777047293527
$overview.Range("F:F").ColumnWidth = 29.9777047293527

$zhcn.Range("C:C").ColumnWidth = 29.9777047293527
$zhcn.Range("P:P").ColumnWidth = 13.7470528738839

$dede.Range("C:C").ColumnWidth = 29.9777047293527
$dede.Range("P:P").ColumnWidth = 13.7470528738839
